$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.03887266666666667
$ws.Range("H2").Value = 0.116618
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1863256666666667
$ws.Range("N2").Value = 0.5589770000000001
$ws.Range("O2").Value = 0.01657678358851065
$ws.Range("P2").Value = 0.01657678358851065
$ws.Range("Q2").Value = 0.007242975531777779
$ws.Range("R2").Value = 0.065186779786
$ws.Range("S2").Value = 0.01657678358851065
$ws.Range("T2").Value = 0.01657678358851065

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.03887266666666667
$ws.Range("H3").Value = 0.116618
$ws.Range("O3").Value = 0.5186672939413604
$ws.Range("P3").Value = 0.5186672939413604
$ws.Range("Q3").Value = 0.2266238501028889
$ws.Range("R3").Value = 2.039614650926
$ws.Range("S3").Value = 0.5186672939413604
$ws.Range("T3").Value = 0.5186672939413604

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.03887266666666667
$ws.Range("H4").Value = 0.116618
$ws.Range("M4").Value = 5.212463666666667
$ws.Range("N4").Value = 15.637391
$ws.Range("O4").Value = 0.4637358003923671
$ws.Range("P4").Value = 0.4637358003923669
$ws.Range("Q4").Value = 0.2026223626264445
$ws.Range("R4").Value = 1.823601263638
$ws.Range("S4").Value = 0.4637358003923671
$ws.Range("T4").Value = 0.4637358003923669

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.03887266666666667
$ws.Range("H5").Value = 0.116618
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.01146633333333333
$ws.Range("N5").Value = 0.034399
$ws.Range("O5").Value = 0.001020122077762015
$ws.Range("P5").Value = 0.001020122077762015
$ws.Range("Q5").Value = 0.0004457269535555556
$ws.Range("R5").Value = 0.004011542582
$ws.Range("S5").Value = 0.001020122077762015
$ws.Range("T5").Value = 0.001020122077762015
